$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Дина"
$ws.Range("B3").Value = "г. Астана"
$ws.Range("B4").Value = "Надо было"
$ws.Range("B5").Value = "12.12.2022 - 24.12.2022"

$ws.Range("B6").Value = 100000.0
$ws.Range("B10").Value = 10000.0
$ws.Range("B11").Value = 60000.0
$ws.Range("B12").Value = 15000.0
$ws.Range("B13").Value = 85000.0
$ws.Range("B14").Value = 15000.0
